$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "subteams" table (Tabelle1) currently spans A1:B10. Add a new row for
# the "StudyDesign" subteam / objective (this grows the table to A1:B11,
# updates the autofilter ref, and appends a new shared string pair).
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

$ws.Range("A11").Value = "StudyDesign"
$ws.Range("B11").Value = "* Investigate the impact of estimands on trial design.`n* Develop approaches and guidance for aligning the design of a trial to its estimands. "

# Match the formatting used by the rest of the table: column A top aligned,
# column B top aligned + wrapped, row height sized for the two wrapped lines.
$ws.Range("A11").VerticalAlignment = -4160
$ws.Range("B11").VerticalAlignment = -4160
$ws.Range("B11").WrapText = $true
$ws.Rows.Item(11).RowHeight = 29

# Move the selection the way the authored workbook ended up (just past the
# newly added row).
$null = $ws.Range("A12").Select()
